$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Benchmark measurement refresh (MSVC 19.42.34436 test run) ---
# Row 5
$ws.Range("B5").Value = 507
$ws.Range("C5").Value = 1176
$ws.Range("D5").Value = 576
$ws.Range("E5").Value = 347
$ws.Range("F5").Value = 349
$ws.Range("G5").Value = 351
$ws.Range("H5").Value = 346
$ws.Range("I5").Value = 1263

# Row 6
$ws.Range("B6").Value = 350
$ws.Range("C6").Value = 347
$ws.Range("D6").Value = 346
$ws.Range("E6").Value = 347
$ws.Range("F6").Value = 346
$ws.Range("G6").Value = 346
$ws.Range("H6").Value = 346
$ws.Range("I6").Value = 348

# Row 7
$ws.Range("B7").Value = 352
$ws.Range("C7").Value = 349
$ws.Range("D7").Value = 347
$ws.Range("E7").Value = 352
$ws.Range("F7").Value = 347
$ws.Range("G7").Value = 347
$ws.Range("H7").Value = 347
$ws.Range("I7").Value = 351

# Row 8
$ws.Range("B8").Value = 165
$ws.Range("C8").Value = 166
$ws.Range("D8").Value = 167
$ws.Range("E8").Value = 168
$ws.Range("F8").Value = 168
$ws.Range("G8").Value = 167
$ws.Range("H8").Value = 166
$ws.Range("I8").Value = 166

# Row 9
$ws.Range("B9").Value = 144
$ws.Range("C9").Value = 146
$ws.Range("D9").Value = 145
$ws.Range("E9").Value = 147
$ws.Range("F9").Value = 146
$ws.Range("G9").Value = 146
$ws.Range("H9").Value = 148
$ws.Range("I9").Value = 146

# Row 10
$ws.Range("B10").Value = 345
$ws.Range("C10").Value = 346
$ws.Range("D10").Value = 348
$ws.Range("E10").Value = 347
$ws.Range("F10").Value = 348
$ws.Range("G10").Value = 349
$ws.Range("H10").Value = 350
$ws.Range("I10").Value = 347

# Row 11
$ws.Range("B11").Value = 348
$ws.Range("C11").Value = 347
$ws.Range("D11").Value = 349
$ws.Range("E11").Value = 347
$ws.Range("F11").Value = 346
$ws.Range("G11").Value = 354
$ws.Range("H11").Value = 346
$ws.Range("I11").Value = 347

# Row 12
$ws.Range("B12").Value = 160
$ws.Range("C12").Value = 160
$ws.Range("D12").Value = 160
$ws.Range("E12").Value = 162
$ws.Range("F12").Value = 159
$ws.Range("G12").Value = 161
$ws.Range("H12").Value = 162
$ws.Range("I12").Value = 161

# Row 13
$ws.Range("B13").Value = 160
$ws.Range("C13").Value = 160
$ws.Range("D13").Value = 159
$ws.Range("E13").Value = 160
$ws.Range("F13").Value = 160
$ws.Range("G13").Value = 159
$ws.Range("H13").Value = 159
$ws.Range("I13").Value = 160

# Row 14
$ws.Range("B14").Value = 144
$ws.Range("C14").Value = 143
$ws.Range("D14").Value = 142
$ws.Range("E14").Value = 147
$ws.Range("F14").Value = 143
$ws.Range("G14").Value = 143
$ws.Range("H14").Value = 143
$ws.Range("I14").Value = 143

# Row 15
$ws.Range("B15").Value = 3274
$ws.Range("C15").Value = 3980
$ws.Range("D15").Value = 16122
$ws.Range("E15").Value = 1727
$ws.Range("F15").Value = 1505
$ws.Range("G15").Value = 1599
$ws.Range("H15").Value = 1937
$ws.Range("I15").Value = 2013

# Row 16
$ws.Range("B16").Value = 532
$ws.Range("C16").Value = 374
$ws.Range("D16").Value = 355
$ws.Range("E16").Value = 390
$ws.Range("F16").Value = 374
$ws.Range("G16").Value = 525
$ws.Range("H16").Value = 370
$ws.Range("I16").Value = 393

# Row 17
$ws.Range("B17").Value = 171
$ws.Range("C17").Value = 396
$ws.Range("D17").Value = 360
$ws.Range("E17").Value = 171
$ws.Range("F17").Value = 168
$ws.Range("G17").Value = 169
$ws.Range("H17").Value = 168
$ws.Range("I17").Value = 372

# Row 18
$ws.Range("B18").Value = 169
$ws.Range("C18").Value = 373
$ws.Range("D18").Value = 358
$ws.Range("E18").Value = 210
$ws.Range("F18").Value = 171
$ws.Range("G18").Value = 173
$ws.Range("H18").Value = 172
$ws.Range("I18").Value = 373

# Row 19
$ws.Range("B19").Value = 162
$ws.Range("C19").Value = 373
$ws.Range("D19").Value = 356
$ws.Range("E19").Value = 179
$ws.Range("F19").Value = 194
$ws.Range("G19").Value = 158
$ws.Range("H19").Value = 158
$ws.Range("I19").Value = 374

# Row 20
$ws.Range("B20").Value = 9172
$ws.Range("C20").Value = 15245
$ws.Range("D20").Value = 4250
$ws.Range("E20").Value = 6778
$ws.Range("F20").Value = 391
$ws.Range("G20").Value = 374
$ws.Range("H20").Value = 1389
$ws.Range("I20").Value = 1269

# Row 21
$ws.Range("B21").Value = 9053
$ws.Range("C21").Value = 15079
$ws.Range("D21").Value = 4338
$ws.Range("E21").Value = 6903
$ws.Range("F21").Value = 361
$ws.Range("G21").Value = 378
$ws.Range("H21").Value = 1411
$ws.Range("I21").Value = 1251

# Row 22
$ws.Range("B22").Value = 167
$ws.Range("C22").Value = 167
$ws.Range("D22").Value = 170
$ws.Range("E22").Value = 169
$ws.Range("F22").Value = 165
$ws.Range("G22").Value = 168
$ws.Range("H22").Value = 178
$ws.Range("I22").Value = 167

# Row 23
$ws.Range("B23").Value = 170
$ws.Range("C23").Value = 179
$ws.Range("D23").Value = 172
$ws.Range("E23").Value = 167
$ws.Range("F23").Value = 169
$ws.Range("G23").Value = 167
$ws.Range("H23").Value = 168
$ws.Range("I23").Value = 169

# Row 24
$ws.Range("B24").Value = 157
$ws.Range("C24").Value = 157
$ws.Range("D24").Value = 157
$ws.Range("E24").Value = 169
$ws.Range("F24").Value = 160
$ws.Range("G24").Value = 161
$ws.Range("H24").Value = 157
$ws.Range("I24").Value = 159

# --- Environment/version labels ---
$ws.Range("A28").Value = "gcc 14.2"
$ws.Range("B28").Value = "Linux Fedora 41"
$ws.Range("A29").Value = "clang 19.1"
$ws.Range("B29").Value = "Linux Fedora 41"
$ws.Range("A30").Value = "Visual Studio 2022,  19.42.34436"
$ws.Range("D30").Value = "10.0.22631.4751"

# --- Restore current selection marker (cosmetic, matches authors last-saved cell) ---
$ws.Range("I45").Select()
